# Update "Maximum Capacity Factor.xlsx"
#   - About!C1: bump the "last updated" date serial 45320 -> 45392 (2024-01-29 -> 2024-04-10)
#   - MCF sheet: set every capacity-factor value that is not already 1 (and not 0) to 1
#   - MCF sheet: move the active cell selection to B17

$wb = $excel.ActiveWorkbook

# --- "About" sheet -----------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- "MCF" sheet --------------------------------------------------------
$wsMcf = $wb.Worksheets.Item("MCF")

# Raise all maximum capacity factors up to 1 (rows 2-18, column B)
foreach ($r in 2..18) {
    $cell = $wsMcf.Cells.Item($r, 2)
    $v = $cell.Value()
    if ($v -gt 0 -and $v -lt 1) {
        $cell.Value = 1
    }
}

# Move the visible selection to B17 to match the saved view state
$wsMcf.Activate()
$wsMcf.Range("B17").Select()
